# Update cryptocurrency price/volume data per latest GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$s = $cell.Style
$cell.Value = '''34.875.00'
$cell.Style = $s
$cell = $ws.Range("E2")
$s = $cell.Style
$cell.Value = '''  -2.31%  '
$cell.Style = $s
$cell = $ws.Range("D3")
$s = $cell.Style
$cell.Value = '''1.808.22'
$cell.Style = $s
$cell = $ws.Range("E3")
$s = $cell.Style
$cell.Value = '''  -3.01%  '
$cell.Style = $s
$cell = $ws.Range("E4")
$s = $cell.Style
$cell.Value = '''  +0.17%  '
$cell.Style = $s
$cell = $ws.Range("D5")
$s = $cell.Style
$cell.Value = '''230.78'
$cell.Style = $s
$cell = $ws.Range("E5")
$s = $cell.Style
$cell.Value = '''  -0.25%  '
$cell.Style = $s
$cell = $ws.Range("E6")
$s = $cell.Style
$cell.Value = '''  -1.36%  '
$cell.Style = $s
$cell = $ws.Range("E7")
$s = $cell.Style
$cell.Value = '''  +0.23%  '
$cell.Style = $s
$cell = $ws.Range("D8")
$s = $cell.Style
$cell.Value = '''39.07'
$cell.Style = $s
$cell = $ws.Range("E8")
$s = $cell.Style
$cell.Value = '''  -8.52%  '
$cell.Style = $s
$cell = $ws.Range("D9")
$s = $cell.Style
$cell.Value = '''0.319'
$cell.Style = $s
$cell = $ws.Range("E9")
$s = $cell.Style
$cell.Value = '''  +2.61%  '
$cell.Style = $s
$cell = $ws.Range("E10")
$s = $cell.Style
$cell.Value = '''  -2.75%  '
$cell.Style = $s
$cell = $ws.Range("E11")
$s = $cell.Style
$cell.Value = '''  -2.14%  '
$cell.Style = $s
$cell = $ws.Range("D12")
$s = $cell.Style
$cell.Value = '''2.069.77'
$cell.Style = $s
$cell = $ws.Range("E12")
$s = $cell.Style
$cell.Value = '''  -3.03%  '
$cell.Style = $s
$cell = $ws.Range("D13")
$s = $cell.Style
$cell.Value = '''1.803.27'
$cell.Style = $s
$cell = $ws.Range("E13")
$s = $cell.Style
$cell.Value = '''  -3.39%  '
$cell.Style = $s
$cell = $ws.Range("E14")
$s = $cell.Style
$cell.Value = '''  -3.28%  '
$cell.Style = $s
$cell = $ws.Range("D15")
$s = $cell.Style
$cell.Value = '''10.85'
$cell.Style = $s
$cell = $ws.Range("E15")
$s = $cell.Style
$cell.Value = '''  -7.31%  '
$cell.Style = $s
$cell = $ws.Range("E16")
$s = $cell.Style
$cell.Value = '''  -4.43%  '
$cell.Style = $s
$cell = $ws.Range("D17")
$s = $cell.Style
$cell.Value = '''34.843.70'
$cell.Style = $s
$cell = $ws.Range("E17")
$s = $cell.Style
$cell.Value = '''  -2.41%  '
$cell.Style = $s
$cell = $ws.Range("D18")
$s = $cell.Style
$cell.Value = '''69.28'
$cell.Style = $s
$cell = $ws.Range("E18")
$s = $cell.Style
$cell.Value = '''  -1.87%  '
$cell.Style = $s
$cell = $ws.Range("D19")
$s = $cell.Style
$cell.Value = '''0.0₃0781'
$cell.Style = $s
$cell = $ws.Range("E19")
$s = $cell.Style
$cell.Value = '''  -3.17%  '
$cell.Style = $s
$cell = $ws.Range("D20")
$s = $cell.Style
$cell.Value = '''238.95'
$cell.Style = $s
$cell = $ws.Range("E20")
$s = $cell.Style
$cell.Value = '''  -4.15%  '
$cell.Style = $s
$cell = $ws.Range("D21")
$s = $cell.Style
$cell.Value = '''11.73'
$cell.Style = $s
$cell = $ws.Range("E21")
$s = $cell.Style
$cell.Value = '''  -4.93%  '
$cell.Style = $s
$cell = $ws.Range("E22")
$s = $cell.Style
$cell.Value = '''  -2.43%  '
$cell.Style = $s
$cell = $ws.Range("E23")
$s = $cell.Style
$cell.Value = '''  +0.22%  '
$cell.Style = $s
$cell = $ws.Range("E24")
$s = $cell.Style
$cell.Value = '''  -0.86%  '
$cell.Style = $s
$cell = $ws.Range("D25")
$s = $cell.Style
$cell.Value = '''173.58'
$cell.Style = $s
$cell = $ws.Range("E25")
$s = $cell.Style
$cell.Value = '''  +1.62%  '
$cell.Style = $s
$cell = $ws.Range("D26")
$s = $cell.Style
$cell.Value = '''7.75'
$cell.Style = $s
$cell = $ws.Range("E26")
$s = $cell.Style
$cell.Value = '''  -3.39%  '
$cell.Style = $s
$cell = $ws.Range("E27")
$s = $cell.Style
$cell.Value = '''  -4.21%  '
$cell.Style = $s
$cell = $ws.Range("E28")
$s = $cell.Style
$cell.Value = '''  -3.40%  '
$cell.Style = $s
$cell = $ws.Range("E29")
$s = $cell.Style
$cell.Value = '''  +5.18%  '
$cell.Style = $s
$cell = $ws.Range("E30")
$s = $cell.Style
$cell.Value = '''  +0.13%  '
$cell.Style = $s
$cell = $ws.Range("E31")
$s = $cell.Style
$cell.Value = '''  +0.40%  '
$cell.Style = $s
$cell = $ws.Range("E32")
$s = $cell.Style
$cell.Value = '''  -0.56%  '
$cell.Style = $s
$cell = $ws.Range("E33")
$s = $cell.Style
$cell.Value = '''  -4.33%  '
$cell.Style = $s
$cell = $ws.Range("D34")
$s = $cell.Style
$cell.Value = '''1.18'
$cell.Style = $s
$cell = $ws.Range("E34")
$s = $cell.Style
$cell.Value = '''  +7.57%  '
$cell.Style = $s
$cell = $ws.Range("E35")
$s = $cell.Style
$cell.Value = '''  -8.12%  '
$cell.Style = $s
$cell = $ws.Range("D36")
$s = $cell.Style
$cell.Value = '''0.683'
$cell.Style = $s
$cell = $ws.Range("E36")
$s = $cell.Style
$cell.Value = '''  -1.34%  '
$cell.Style = $s
$cell = $ws.Range("D37")
$s = $cell.Style
$cell.Value = '''90.66'
$cell.Style = $s
$cell = $ws.Range("E37")
$s = $cell.Style
$cell.Value = '''  -10.74%  '
$cell.Style = $s
$cell = $ws.Range("E38")
$s = $cell.Style
$cell.Value = '''  +6.15%  '
$cell.Style = $s
$cell = $ws.Range("D39")
$s = $cell.Style
$cell.Value = '''1.308.25'
$cell.Style = $s
$cell = $ws.Range("E39")
$s = $cell.Style
$cell.Value = '''  -4.65%  '
$cell.Style = $s
$cell = $ws.Range("E40")
$s = $cell.Style
$cell.Value = '''  -3.06%  '
$cell.Style = $s
$cell = $ws.Range("E41")
$s = $cell.Style
$cell.Value = '''  -0.69%  '
$cell.Style = $s
$cell = $ws.Range("E42")
$s = $cell.Style
$cell.Value = '''  -6.19%  '
$cell.Style = $s
$cell = $ws.Range("D43")
$s = $cell.Style
$cell.Value = '''14.12'
$cell.Style = $s
$cell = $ws.Range("E43")
$s = $cell.Style
$cell.Value = '''  -5.82%  '
$cell.Style = $s
$cell = $ws.Range("D44")
$s = $cell.Style
$cell.Value = '''2.19'
$cell.Style = $s
$cell = $ws.Range("E44")
$s = $cell.Style
$cell.Value = '''  -12.81%  '
$cell.Style = $s
$cell = $ws.Range("D45")
$s = $cell.Style
$cell.Value = '''2.70'
$cell.Style = $s
$cell = $ws.Range("B46")
$s = $cell.Style
$cell.Value = '''Kaspa'
$cell.Style = $s
$cell = $ws.Range("C46")
$s = $cell.Style
$cell.Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell.Style = $s
$cell = $ws.Range("D46")
$s = $cell.Style
$cell.Value = '''0.0510'
$cell.Style = $s
$cell = $ws.Range("E46")
$s = $cell.Style
$cell.Value = '''  -1.93%  '
$cell.Style = $s
$cell = $ws.Range("B47")
$s = $cell.Style
$cell.Value = '''FraxShare'
$cell.Style = $s
$cell = $ws.Range("C47")
$s = $cell.Style
$cell.Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell.Style = $s
$cell = $ws.Range("D47")
$s = $cell.Style
$cell.Value = '''6.10'
$cell.Style = $s
$cell = $ws.Range("E47")
$s = $cell.Style
$cell.Value = '''  -3.20%  '
$cell.Style = $s
$cell = $ws.Range("D48")
$s = $cell.Style
$cell.Value = '''1.996.11'
$cell.Style = $s
$cell = $ws.Range("E48")
$s = $cell.Style
$cell.Value = '''  -1.87%  '
$cell.Style = $s
$cell = $ws.Range("B49")
$s = $cell.Style
$cell.Value = '''Cronos'
$cell.Style = $s
$cell = $ws.Range("C49")
$s = $cell.Style
$cell.Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell.Style = $s
$cell = $ws.Range("D49")
$s = $cell.Style
$cell.Value = '''0.0673'
$cell.Style = $s
$cell = $ws.Range("E49")
$s = $cell.Style
$cell.Value = '''  +7.28%  '
$cell.Style = $s
$cell = $ws.Range("B50")
$s = $cell.Style
$cell.Value = '''PaxDollar'
$cell.Style = $s
$cell = $ws.Range("C50")
$s = $cell.Style
$cell.Value = '''https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell.Style = $s
$cell = $ws.Range("D50")
$s = $cell.Style
$cell.Value = '''1.01'
$cell.Style = $s
$cell = $ws.Range("E50")
$s = $cell.Style
$cell.Value = '''  +0.21%  '
$cell.Style = $s
$cell = $ws.Range("D51")
$s = $cell.Style
$cell.Value = '''98.63'
$cell.Style = $s
$cell = $ws.Range("E51")
$s = $cell.Style
$cell.Value = '''  -6.15%  '
$cell.Style = $s
